$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column "State CV" values for rows 6-10
$ws.Range("C6").Value = 0.99
$ws.Range("C7").Value = 0.9
$ws.Range("C8").Value = 0.89
$ws.Range("C9").Value = 1.03
$ws.Range("C10").Value = 1.01

# Add "Full model" values for ann (row 11) and keras (row 12)
$ws.Range("B11").Value = 0.86
$ws.Range("B12").Value = 0.89

# Add new note row at 15
$ws.Range("A15").Value = "ann probably need normalization to work better"

# Update selection to match target
$ws.Range("B13").Select()
